# "Add files via upload" - update World_pop_proj.xlsx
#  - rename shared-string column headers: "census" -> "US_Census", "un" -> "UN"
#  - move the active selection on the worksheet to F15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "US_Census"
$ws.Range("C1").Value = "UN"

$ws.Range("F15").Select() | Out-Null
